# DTranslate.pptx — "Important Terms Check" -> "Key Word Check" rename
#
# Renames the "Important Terms Check" / "Important terms check" MAST-step
# label to "Key Word Check" / "Key Word check" everywhere it appears in the
# deck, and nudges the repositioned label box on the MAST-checking-steps
# slide to match the shorter text's new (re-centered) placement.

$p = $ppt.ActivePresentation

# EMU-exact point value for a target-EMU shape coordinate: the COM bridge
# stores shape geometry as single-precision floats and floors on EMU
# round-trip, so biasing by half an EMU before converting keeps the
# save at the exact integer EMU we want instead of landing one unit short.
function PtForEmu([double]$emu) {
    return ($emu + 0.5) / 12700.0
}

# ---------------------------------------------------------------------
# Slide 12 — bullet list, step 7 bullet text
# ---------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$sh12 = $s12.Shapes.Item(4)   # "Content Placeholder 2"
$run12 = $sh12.TextFrame.TextRange.Paragraphs(7).Runs(1)
$run12.Text = "Key Word check: Check key terms to ensure they are present in the draft and translated clearly and consistently."

# ---------------------------------------------------------------------
# Slide 13 — bullet list, step 7 bullet text
# ---------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$sh13 = $s13.Shapes.Item(4)   # "Content Placeholder 2"
$run13 = $sh13.TextFrame.TextRange.Paragraphs(7).Runs(1)
$run13.Text = "Key Word check: Check key terms to ensure they are present in the draft and translated clearly and consistently."

# ---------------------------------------------------------------------
# Slide 14 — MAST diagram label textbox ("TextBox 29")
# spAutoFit recalculates the box height for the shorter caption, so the
# original height is restored afterwards to keep the shape geometry
# untouched (matches the source diff, which only changes the text).
# ---------------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$sh14 = $s14.Shapes.Item(17)  # "TextBox 29"
$origHeight14 = $sh14.Height
$run14 = $sh14.TextFrame.TextRange.Paragraphs(1).Runs(1)
$run14.Text = "Key Word Check"
$sh14.Height = $origHeight14

# ---------------------------------------------------------------------
# Slide 23 — grouped label textbox ("TextBox 14") under the 3rd icon
# Text shrinks and the box is re-centered under its picture, so position
# and width move to the new (smaller) box recorded in the target diff.
# ---------------------------------------------------------------------
$s23 = $p.Slides.Item(23)
$grp23 = $s23.Shapes.Item(7)        # "Group 36"
$sh23 = $grp23.GroupItems.Item(3)   # "TextBox 14"
$run23 = $sh23.TextFrame.TextRange.Paragraphs(1).Runs(1)
$run23.Text = "Key Word Check"
$sh23.Left = PtForEmu 6877956
$sh23.Width = PtForEmu 2066591

# ---------------------------------------------------------------------
# Slide 34 — slide title
# ---------------------------------------------------------------------
$s34 = $p.Slides.Item(34)
$sh34 = $s34.Shapes.Item(3)   # "Title 1"
$run34 = $sh34.TextFrame.TextRange.Paragraphs(1).Runs(1)
$run34.Text = "MAST Step 7: Key Word Check "
